$wb = $excel.ActiveWorkbook

# Rename sheet "Wong3" to "Euclid"
$ws = $wb.Worksheets.Item("Wong3")
$ws.Name = "Euclid"

# Update the SBFL data (columns C and D) for the changed rows
$updates = @(
    @{Row=2;  C=229; D=51.11607142857143},
    @{Row=3;  C=250; D=55.80357142857143},
    @{Row=4;  C=439; D=97.99107142857143},
    @{Row=5;  C=420; D=93.75},
    @{Row=6;  C=395; D=88.16964285714286},
    @{Row=7;  C=148; D=33.03571428571428},
    @{Row=8;  C=322; D=71.875},
    @{Row=9;  C=196; D=43.75},
    @{Row=10; C=336; D=75},
    @{Row=11; C=224; D=50},
    @{Row=12; C=331; D=73.88392857142857},
    @{Row=14; C=176; D=39.28571428571428},
    @{Row=15; C=348; D=77.67857142857143},
    @{Row=16; C=399; D=89.0625},
    @{Row=17; C=130; D=29.01785714285715},
    @{Row=18; C=110; D=24.55357142857143},
    @{Row=21; C=395; D=88.16964285714286}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 3).Value = $u.C
    $ws.Cells.Item($u.Row, 4).Value = $u.D
}
